$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materiales para TABLEROS")

# --- Quantity (column E) updates ---
$ws.Range("E10").Value = 450
$ws.Range("E12").Value = 60
$ws.Range("E17").Value = 20
$ws.Range("E18").Value = 20
$ws.Range("E19").Value = 50
$ws.Range("E20").Value = 50
$ws.Range("E21").Value = 200
$ws.Range("E22").Value = 6
$ws.Range("E23").Value = 15
$ws.Range("E24").Value = 0
$ws.Range("E26").Value = 300
$ws.Range("E28").Value = 1

# --- Clear stray note text in B42 ("agregar bornes de hacia variador") ---
$ws.Range("B42").ClearContents()

# --- View changes: zoom + frozen pane top-left cell + selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("B11").Select()

$wb.Save()
